$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.147441625595093
$ws.Range("B1").Value = 2.402315139770508
$ws.Range("C1").Value = 5.151350975036621
$ws.Range("D1").Value = 2.213490962982178
$ws.Range("E1").Value = 1.24798309803009
